# Update the "Generate Report for Handback" timestamps.
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first row
$wsOverview.Range("G2").Value = "2016-08-16 21:01:22"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for the first row
$wsZhCn.Range("H2").Value = "2016-08-16 21:01:17"
$wsZhCn.Range("K2").Value = "2016-08-16 21:01:34"

# de-de sheet: "Correspond Handoff Datetime" (mirrors Overview's value) and
# "Correspond Handback DateTime" for the first row
$wsDeDe.Range("H2").Value = "2016-08-16 21:01:22"
$wsDeDe.Range("K2").Value = "2016-08-16 21:01:41"
